# Auto-generated script to apply F-column ('想去人数' / interest count) updates
# across all four worksheets, matching the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1740
$ws.Range("F4").Value = 150
$ws.Range("F5").Value = 402
$ws.Range("F6").Value = 785
$ws.Range("F7").Value = 226
$ws.Range("F8").Value = 1124
$ws.Range("F9").Value = 305
$ws.Range("F11").Value = 862
$ws.Range("F12").Value = 657
$ws.Range("F13").Value = 178
$ws.Range("F17").Value = 163
$ws.Range("F18").Value = 2867
$ws.Range("F19").Value = 2602
$ws.Range("F23").Value = 311
$ws.Range("F26").Value = 5210
$ws.Range("F29").Value = 16
$ws.Range("F31").Value = 276
$ws.Range("F32").Value = 1059

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1075
$ws.Range("F5").Value = 1075
$ws.Range("F7").Value = 17
$ws.Range("F8").Value = 231
$ws.Range("F11").Value = 13
$ws.Range("F15").Value = 599
$ws.Range("F26").Value = 271
$ws.Range("F27").Value = 3860
$ws.Range("F29").Value = 5
$ws.Range("F31").Value = 196
$ws.Range("F32").Value = 39
$ws.Range("F34").Value = 151

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1012
$ws.Range("F9").Value = 1272
$ws.Range("F10").Value = 342
$ws.Range("F11").Value = 93

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1740
$ws.Range("F6").Value = 1012
$ws.Range("F7").Value = 1272
$ws.Range("F8").Value = 342
$ws.Range("F9").Value = 93
$ws.Range("F10").Value = 150
$ws.Range("F11").Value = 402
$ws.Range("F12").Value = 785
$ws.Range("F13").Value = 226
$ws.Range("F15").Value = 1124
$ws.Range("F16").Value = 305
$ws.Range("F17").Value = 657
$ws.Range("F18").Value = 1075
$ws.Range("F19").Value = 178
$ws.Range("F21").Value = 17
$ws.Range("F22").Value = 163
$ws.Range("F23").Value = 2867
$ws.Range("F24").Value = 2602
$ws.Range("F25").Value = 231
$ws.Range("F26").Value = 311
$ws.Range("F28").Value = 13
$ws.Range("F31").Value = 5210
$ws.Range("F34").Value = 599
$ws.Range("F35").Value = 599
$ws.Range("F36").Value = 16
$ws.Range("F38").Value = 276
$ws.Range("F43").Value = 271
$ws.Range("F44").Value = 1059
$ws.Range("F45").Value = 196
$ws.Range("F46").Value = 39
$ws.Range("F48").Value = 151

$wb.Save()
